# Weekly update: a new observation (week) is inserted at row 6 of the
# "Arveja Verde" sheet, pushing all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 6..70 down to 7..71, carrying formatting (e.g. the date
# style on column D) the same way Excel's own Insert does.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with this week's record.
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").Value = 44630
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = 100112022
$ws.Range("G6").Value = "Arveja Verde"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 60
$ws.Range("K6").Value = 24000
$ws.Range("L6").Value = 25000
$ws.Range("M6").Value = 24500
$ws.Range("N6").Value = "$/saco 25 kilos"
$ws.Range("O6").Value = "Provincia de Diguillín"
$ws.Range("P6").Value = 980
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = "Hortaliza"
